$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same style (bold, bordered, centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header row (row 1) - new columns I0 and IF
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-16
$data = @(
    @(11, 11),
    @(6, 7),
    @(8, 8),
    @(4, 5),
    @(7, 7),
    @(9, 9),
    @(6, 7),
    @(7, 7),
    @(6, 7),
    @(8, 8),
    @(5, 5),
    @(9, 9),
    @(6, 6),
    @(9, 9),
    @(7, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
